$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy the formatting for columns A:O from row 12 down to the new row 13,
#    reusing the existing styles (s=6, s=3, s=3, default, default, default,
#    default, default, s=5 x7) exactly as used on row 12.
$ws.Range("A12:O12").Copy()
$ws.Range("A13:O13").PasteSpecial(-4122)  # xlPasteFormats

# 2) Build the new cell style for P13:Q13 (numFmtId 49 "@" + centered,
#    no wrap) by formatting a single cell first (so the style engine
#    resolves cleanly to one new xf), then propagate it via copy/paste
#    of formats only, to avoid creating stray intermediate styles.
$p13 = $ws.Cells.Item(13, 16)
$p13.NumberFormat = "@"
$p13.WrapText = $false
$p13.HorizontalAlignment = -4108  # xlCenter

$ws.Cells.Item(13, 16).Copy()
$ws.Cells.Item(13, 17).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3) Fill in the values for the new row.
$ws.Range("A13").Value = "namespace_style_invalid"
$ws.Range("B13").Value = "A namespace must have a valid style."
$ws.Range("C13").Value = "Invalid value"
$ws.Range("D13").Value = "Namespace"
$ws.Range("E13").Value = "Style"
$ws.Range("F13").Value = "All"
$ws.Range("G13").Value = "spreadsheet"
$ws.Range("H13").Value = "error"

# 4) Grow the table (ListObject) so it covers the new row; this also keeps
#    the autoFilter range, dimension, etc. in sync.
$tbl = $ws.ListObjects.Item(1)
$newRange = $ws.Range("A1:U13")
$tbl.Resize($newRange)

# 5) Match the cursor/selection move down to the newly added row.
[void]$ws.Range("A13").Select()

Write-Output "edit complete"
